# This script reproduces the diff:
#  - 2 new rows are inserted above the original header row.
#  - New row 1 becomes a numeric index row (0..13), taking on the bold/
#    centered/thin-bordered header style that used to live on row 1.
#  - New row 2 is an almost-empty row, with only "Washer" in column E.
#  - New row 3 holds the original text header labels (Lg., Threading, ...)
#    but WITHOUT the special header styling (plain formatting, like a
#    normal data row).
#  - The rest of the data (old rows 2-35) shifts down to rows 4-37
#    unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the header labels that currently live in row 1 before we
# shift everything down.
$headerValues = @(
    "Lg.",
    "Threading",
    "HeadDia.",
    "Head Ht.",
    "OD",
    "Thick.",
    "DriveSize",
    "TensileStrength, psi",
    "Specifications Met",
    "Pkg.Qty.",
    "",
    "Pkg.",
    "thread_size",
    "material_surface"
)

# Insert two new blank rows above row 1; this pushes the current
# contents of row 1 (headers, with their bold/centered/bordered style)
# down to row 3, and the rest of the table down accordingly
# (old row 2 -> row 4, ..., old row 35 -> row 37). The two freshly
# inserted rows (1 and 2) come out with plain/default formatting.
$ws.Rows("1:2").Insert()

# --- Row 3: plain (unstyled) header labels -------------------------------
# Row 3 inherited the old bold/centered/bordered header formatting from
# the Insert() above. Strip it back to the plain/default look used by
# normal data rows, then (re)write the header text.
$headerRange = $ws.Range("A3:N3")
$headerRange.ClearFormats()

for ($i = 0; $i -lt $headerValues.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(3, $col).Value = $headerValues[$i]
}

# --- Row 2: the "Washer" marker row --------------------------------------
for ($col = 1; $col -le 14; $col++) {
    $ws.Cells.Item(2, $col).Value = ""
}
$ws.Cells.Item(2, 5).Value = "Washer"

# --- Row 1: numeric column-index row, gets the header styling -----------
for ($col = 1; $col -le 14; $col++) {
    $ws.Cells.Item(1, $col).Value = $col - 1
}

$indexRange = $ws.Range("A1:N1")
$indexRange.Font.Bold = $true
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
$indexRange.Borders.LineStyle = 1
$indexRange.Borders.Weight = 2
